$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9724684953689575
$ws.Range("B1").Value = 1.749895930290222
$ws.Range("C1").Value = 4.794293403625488
$ws.Range("D1").Value = 1.356752753257751
$ws.Range("E1").Value = 1.23593807220459
